$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 192-193, pushing the existing rows 192.. down to 194..
$ws.Range("A192:A193").EntireRow.Insert()

# New row 192 (new price record dated 45209 / 2023-10-10)
$ws.Range("A192").Value = 10
$ws.Range("B192").Value = "Vega Modelo de Temuco"
$ws.Range("C192").Value = "La Araucanía"
$ws.Range("D192").Value = 45209
$ws.Range("E192").Value = 9
$ws.Range("F192").Value = 100112031
$ws.Range("G192").Value = "Poroto verde"
$ws.Range("H192").Value = "Sin especificar"
$ws.Range("I192").Value = "Primera"
$ws.Range("J192").Value = 80
$ws.Range("K192").Value = 1600
$ws.Range("L192").Value = 1600
$ws.Range("M192").Value = 1600
$ws.Range("N192").Value = "$/kilo"
$ws.Range("O192").Value = "Provincia de Limarí"
$ws.Range("P192").Value = 1600
$ws.Range("Q192").Value = 1
$ws.Range("R192").Value = "Hortaliza"

# New row 193 (new price record dated 45209 / 2023-10-10)
$ws.Range("A193").Value = 10
$ws.Range("B193").Value = "Vega Modelo de Temuco"
$ws.Range("C193").Value = "La Araucanía"
$ws.Range("D193").Value = 45209
$ws.Range("E193").Value = 9
$ws.Range("F193").Value = 100112031
$ws.Range("G193").Value = "Poroto verde"
$ws.Range("H193").Value = "Sin especificar"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 12
$ws.Range("K193").Value = 40000
$ws.Range("L193").Value = 40000
$ws.Range("M193").Value = 40000
$ws.Range("N193").Value = "$/malla 25 kilos"
$ws.Range("O193").Value = "Perú"
$ws.Range("P193").Value = 1600
$ws.Range("Q193").Value = 25
$ws.Range("R193").Value = "Hortaliza"
